# Updates cryptos list: Price (D) and Volume(1h) (E) columns for rows 2-51.
# Some "Price" values look like plain numbers (e.g. "608.44"); a direct
# $range.Value = "608.44" assignment would let Excel auto-convert them to a
# numeric cell (losing the exact decimal text / introducing FP noise). The
# source workbook stores them as literal text, so for those we build the
# text in a scratch cell via a text-formula, copy it, and PasteSpecial just
# the values into the target cell - this keeps the destination a plain text
# cell (no number format change) while still landing exact text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")

function Set-TextValue($cellRef, $text) {
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

$ws.Range('D2').Value = '68.343.38'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '2.716.17'
$ws.Range('E3').Value = '  +2.65%  '
$ws.Range('E4').Value = '  -0.19%  '
Set-TextValue 'D5' '608.44'
$ws.Range('E5').Value = '  +1.97%  '
Set-TextValue 'D6' '167.13'
$ws.Range('E6').Value = '  +5.30%  '
$ws.Range('E7').Value = '  +0.02%  '
Set-TextValue 'D8' '0.554'
$ws.Range('E8').Value = '  +2.73%  '
$ws.Range('D9').Value = '2.715.12'
$ws.Range('E9').Value = '  +2.62%  '
$ws.Range('E10').Value = '  +1.98%  '
Set-TextValue 'D11' '0.365'
$ws.Range('E11').Value = '  +4.27%  '
$ws.Range('E12').Value = '  +0.62%  '
$ws.Range('E13').Value = '  +0.64%  '
Set-TextValue 'D14' '28.55'
$ws.Range('E14').Value = '  +2.32%  '
$ws.Range('D15').Value = '3.217.81'
$ws.Range('E15').Value = '  +2.75%  '
Set-TextValue 'D16' '0.0000188'
$ws.Range('E16').Value = '  +0.94%  '
$ws.Range('D17').Value = '68.274.80'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').Value = '2.651.98'
$ws.Range('E18').Value = '  -0.87%  '
$ws.Range('E19').Value = '  +4.34%  '
Set-TextValue 'D20' '371.71'
$ws.Range('E20').Value = '  +2.25%  '
$ws.Range('E21').Value = '  +2.72%  '
Set-TextValue 'D22' '4.50'
$ws.Range('E22').Value = '  +2.40%  '
Set-TextValue 'D23' '4.98'
$ws.Range('E23').Value = '  +4.63%  '
$ws.Range('E24').Value = '  +1.18%  '
Set-TextValue 'D25' '73.00'
$ws.Range('E25').Value = '  -1.90%  '
$ws.Range('E26').Value = '  +0.04%  '
Set-TextValue 'D27' '10.14'
$ws.Range('E27').Value = '  +4.36%  '
$ws.Range('D28').Value = '2.866.40'
$ws.Range('E28').Value = '  +3.17%  '
$ws.Range('E29').Value = '  +1.53%  '
Set-TextValue 'D30' '583.72'
$ws.Range('E30').Value = '  +3.92%  '
$ws.Range('E31').Value = '  +0.09%  '
Set-TextValue 'D32' '8.21'
$ws.Range('E32').Value = '  +2.13%  '
$ws.Range('E33').Value = '  +2.95%  '
$ws.Range('E34').Value = '  +6.88%  '
$ws.Range('E35').Value = '  +1.87%  '
$ws.Range('E36').Value = '  -3.80%  '
Set-TextValue 'D37' '0.998'
$ws.Range('E37').Value = '  -0.13%  '
Set-TextValue 'D38' '162.89'
$ws.Range('E38').Value = '  +1.98%  '
Set-TextValue 'D39' '19.86'
$ws.Range('E39').Value = '  +1.24%  '
$ws.Range('E40').Value = '  +2.33%  '
Set-TextValue 'D41' '1.88'
$ws.Range('E41').Value = '  +0.80%  '
Set-TextValue 'D42' '5.42'
$ws.Range('E42').Value = '  +2.10%  '
Set-TextValue 'D43' '17.97'
$ws.Range('E43').Value = '  +0.92%  '
$ws.Range('E44').Value = '  +1.05%  '
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('D46').Value = '0.0₆0311'
$ws.Range('E46').Value = '  -2.83%  '
Set-TextValue 'D47' '40.89'
$ws.Range('E47').Value = '  +1.42%  '
Set-TextValue 'D48' '0.597'
$ws.Range('E48').Value = '  +4.17%  '
Set-TextValue 'D49' '155.07'
$ws.Range('E49').Value = '  -1.85%  '
Set-TextValue 'D50' '3.91'
$ws.Range('E50').Value = '  +2.93%  '
$ws.Range('E51').Value = '  +5.20%  '

$scratch.ClearContents()
